$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.11928112952215741
$ws.Range("A2").Value = -0.0059999999401334492
$ws.Range("A3").Value = -0.0039999999517625895
$ws.Range("A4").Value = -0.0079999999096003194
$ws.Range("A5").Value = -0.0029999999562040358
$ws.Range("A6").Value = -0.0019999999595974316
$ws.Range("A7").Value = -0.0099999998807542845
$ws.Range("A8").Value = -0.0099999998812556612
$ws.Range("A9").Value = -0.0019999999621385101
$ws.Range("A10").Value = -0.0019999999648252498
$ws.Range("A11").Value = -0.0029999999551089118
$ws.Range("A12").Value = -0.003499999950824062
$ws.Range("A13").Value = -0.0034999999570457518
$ws.Range("A14").Value = 0.024160878916408279
$ws.Range("A15").Value = -0.00099999998448119243
$ws.Range("A16").Value = -0.001999999975369704
$ws.Range("A17").Value = -0.0019999999772393195
$ws.Range("A18").Value = 0.009920030712113892
$ws.Range("A19").Value = -0.0039999999596487257
$ws.Range("A20").Value = -0.0039999999564805933
$ws.Range("A21").Value = -0.0039999999559974242
$ws.Range("A22").Value = -0.0039999999556137311
$ws.Range("A23").Value = -0.0049999999387413396
$ws.Range("A24").Value = -0.019999999784807265
$ws.Range("A25").Value = -0.019999999781821209
$ws.Range("A26").Value = -0.0024999999537982376
$ws.Range("A27").Value = -0.002499999952660481
$ws.Range("A28").Value = -0.001999999952575493
$ws.Range("A29").Value = -0.0069999999002066104
$ws.Range("A30").Value = -0.059999999378605562
$ws.Range("A31").Value = -0.0069999999080092579
$ws.Range("A32").Value = -0.0099999998814102042
$ws.Range("A33").Value = 0.065483840255774695
